$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the sub-item rows so each includes its section name as a prefix
# (e.g. "     New nominations" -> "     Civilian, New nominations").

# Civilian
$ws.Range("A7").Value  = "     Civilian, New nominations"
$ws.Range("A8").Value  = "     Civilian, Carryover nominations"
$ws.Range("A9").Value  = "     Civilian, Confirmed "
$ws.Range("A10").Value = "     Civilian, Unconfirmed "
$ws.Range("A11").Value = "     Civilian, Withdrawn "
$ws.Range("A12").Value = "     Civilian, Returned to White House "

# Other Civilian
$ws.Range("A14").Value = "     Other Civilian, New nominations"
$ws.Range("A15").Value = "     Other Civilian, Carryover nominations"
$ws.Range("A16").Value = "     Other Civilian, Confirmed "
$ws.Range("A17").Value = "     Other Civilian, Returned to White House "

# Air Force
$ws.Range("A19").Value = "     Air Force, New nominations"
$ws.Range("A20").Value = "     Air Force, Confirmed "
$ws.Range("A21").Value = "     Air Force, Returned to White House "

# Army
$ws.Range("A23").Value = "     Army, New nominations"
$ws.Range("A24").Value = "     Army, Carryover nominations"
$ws.Range("A25").Value = "     Army, Confirmed "
$ws.Range("A26").Value = "     Army, Withdrawn "
$ws.Range("A27").Value = "     Army, Returned to White House "

# Navy
$ws.Range("A29").Value = "     Navy, New nominations"
$ws.Range("A30").Value = "     Navy, Carryover nominations"
$ws.Range("A31").Value = "     Navy, Confirmed "
$ws.Range("A32").Value = "     Navy, Returned to White House "

# Marine Corps
$ws.Range("A34").Value = "     Marine Corps, New nominations"
$ws.Range("A35").Value = "     Marine Corps, Confirmed "
$ws.Range("A36").Value = "     Marine Corps, Returned to White House "

# Space Force
$ws.Range("A38").Value = "     Space Force, New nominations"
$ws.Range("A39").Value = "     Space Force, Confirmed "

# Turn the "Summary" section header into a labelled total row holding the
# overall new-nominations count, and relabel the carried-over total.
$ws.Range("A40").Value = "Total new nominations"
$ws.Range("B40").Value = 21854
$ws.Range("B40").NumberFormat = $ws.Range("B42").NumberFormat

$ws.Range("A41").Value = "Total carryover nominations"

# The old "Total nominations received this Session" row is no longer
# needed now that its value lives on the "Total new nominations" row
# above; remove it and let the remaining total rows shift up.
$ws.Rows(42).Delete()
